$wb = $excel.ActiveWorkbook
$t = $wb.Theme
$tc = $t.ThemeColorScheme
$tc.Item(2).RGB = 16777215
